$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: amount, currency, account_code, description, value_date, business_unit
$ws.Cells.Item(5, 1).Value = 40500
$ws.Cells.Item(5, 2).Value = "BRL"
$ws.Cells.Item(5, 3).Value = "ACC-BU_C-0004"
$ws.Cells.Item(5, 4).Value = "Sample closure line 4 for BU_C"
$ws.Range("E5").NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = "2025-02-15"
$ws.Range("E5").Style = "Normal"
$ws.Cells.Item(5, 6).Value = "BU_C"

# Row 6: amount, currency, account_code, description, value_date, business_unit
$ws.Cells.Item(6, 1).Value = 50500
$ws.Cells.Item(6, 2).Value = "BRL"
$ws.Cells.Item(6, 3).Value = "ACC-BU_C-0005"
$ws.Cells.Item(6, 4).Value = "Sample closure line 5 for BU_C"
$ws.Range("E6").NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = "2025-02-15"
$ws.Range("E6").Style = "Normal"
$ws.Cells.Item(6, 6).Value = "BU_C"
